$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths: both columns A and B become the same width
# (closest achievable value to the target 15.42578125 given this engine's
# internal 1/6-character-width quantization for ColumnWidth).
$ws.Columns.Item(1).ColumnWidth = 14.666666666666668
$ws.Columns.Item(2).ColumnWidth = 14.666666666666668

$values = @(
    @(-0.15874635842922658, 0.15806139970295874),
    @(-0.051127580700061515, 0.04971141197657758),
    @(0.093770658815838459, -0.094452276543258051),
    @(-0.18953586484436613, 0.1882296455530259),
    @(-0.1822296464254638, 0.17957392393540594),
    @(-0.078665239442041734, 0.078577339310698324),
    @(-0.058577340366854358, 0.058397577209083451),
    @(-0.038397578273347222, 0.038289409933624619),
    @(-0.03228941086095638, 0.032212045571005099),
    @(-0.026212046505250441, 0.026211197665574559),
    @(-0.02171119858512327, 0.021702339526658676),
    @(-0.015702340462431685, 0.015677969918904466),
    @(-0.009677970859348406, 0.0096745772560824506),
    @(0.0023254217403625432, -0.0023300782724282598),
    @(0.0083300773312133813, -0.0083408555889779024),
    @(-0.015026978608466024, 0.015003946833059612),
    @(-0.009003947777176613, 0.008999999023942884),
    @(-0.059632230398481312, 0.059581055182412967),
    @(-0.027096231407010674, 0.027013079939922857),
    @(-0.018013080838539253, 0.018004206377439047),
    @(-0.0090042072770861736, 0.008999999099604139),
    @(-0.093937084652933933, 0.093628041692385366),
    @(-0.084628042595770836, 0.084125327326252552),
    @(-0.042125328586860356, 0.041999998732640265),
    @(-0.094908685401208004, 0.094667762933386967),
    @(-0.088667763843627512, 0.088356213009426199),
    @(-0.082356213924460242, 0.081281866556329962),
    @(-0.075281867492562604, 0.0745377324068075),
    @(-0.062537733418958297, 0.062169838871502847),
    @(-0.042169839973326351, 0.042018727160978209),
    @(-0.02701872822042084, 0.027000469855567033),
    @(-0.0060004709785124177, 0.0059999990322632257)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}
